# Revised PM calculation code
#
# The Emissions column (F) previously computed a per-vehicle-type daily
# emission figure (E*D). The revision scales each vehicle category's
# emission factor by its annual/fleet multiplier (vehicle-km-equivalent
# factors), which cascades into the Scenario columns (H, I, J) and the
# Share column (G) through already-existing formulas that reference F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (4W1): E9*D9 -> E9*D9*10000
$ws.Range("F9").Formula = "=E9*D9*10000"
$ws.Range("F9").Style = "Normal"

# --- Row 10 (4W2): E10*D10 -> E10*D10*15000
$ws.Range("F10").Formula = "=E10*D10*15000"
$ws.Range("F10").Style = "Normal"

# --- Row 11 (2W): E11*D11 -> E11*D11*10000 (no longer a shared formula)
$ws.Range("F11").Formula = "=E11*D11*10000"
$ws.Range("F11").Style = "Normal"

# --- Row 12 (3W2S): E12*D12 -> E12*D12*20000
$ws.Range("F12").Formula = "=E12*D12*20000"
$ws.Range("F12").Style = "Normal"

# --- Row 13 (3W4S): E13*D13 -> E13*D13*150*365
$ws.Range("F13").Formula = "=E13*D13*150*365"
$ws.Range("F13").Style = "Normal"

# --- Row 14 (TAXI): E14*D14 -> E14*D14*25*2*365
$ws.Range("F14").Formula = "=E14*D14*25*2*365"
$ws.Range("F14").Style = "Normal"

# --- Row 15 (BUS): E15*D15 -> E15*D15*25*2*365
$ws.Range("F15").Formula = "=E15*D15*25*2*365"
$ws.Range("F15").Style = "Normal"

# --- Row 16 (HDT): E16*D16 -> E16*D16*25*2*365
$ws.Range("F16").Formula = "=E16*D16*25*2*365"
$ws.Range("F16").Style = "Normal"

# The much larger Emissions/Scenario values no longer fit the old column
# widths, so widen column F and the Scenario columns (H:J) to match.
$ws.Range("F1").ColumnWidth = 11.166666666666666
$ws.Range("H1:J1").ColumnWidth = 10.166666666666666
